$p = $ppt.ActivePresentation

# --- Slide 13 : "Instruction conditionnelle « if »" ---
$s13 = $p.Slides.Item(13)

# 1) Remove the stray "Rectangle 21" textbox ("Définir et affecter une valeur à des variables")
for ($i = $s13.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s13.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 21") {
        $shp.Delete()
    }
}

# 2) Reposition "Rectangle 7" (the "si (condition) alors ... else ..." box) lower on the slide
for ($i = 1; $i -le $s13.Shapes.Count; $i++) {
    $shp = $s13.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 7") {
        $shp.Left = 375.7507172614173
        $shp.Top = 222.91858677716536
    }
}

# 3) Fix the "booleén" typo -> "booléen" in "Rectangle 9", merging it into the previous run
for ($i = 1; $i -le $s13.Shapes.Count; $i++) {
    $shp = $s13.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 9") {
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf(" (vrai)")
        if ($idx -ge 0) {
            $len = $full.Length - $idx
            $sub = $tr.Characters($idx + 1, $len)
            $sub.Text = " (vrai) et False (faux) et d’autres expressions entre booléen"
        }
    }
}

# --- Slide 9 : "L'indentation" ---
$s9 = $p.Slides.Item(9)
for ($i = 1; $i -le $s9.Shapes.Count; $i++) {
    $grp = $s9.Shapes.Item($i)
    if ($grp.Name -eq "Groupe 9") {
        for ($j = 1; $j -le $grp.GroupItems.Count; $j++) {
            $shp = $grp.GroupItems.Item($j)
            if ($shp.Name -eq "Rectangle 16") {
                $shp.TextFrame.TextRange.Text = "Le retour à une indentation précédente met fin au « if »"
            }
        }
    }
}
